$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo: C2 should read "FREQ[23:0] hex" (matching the DEC2HEX formula below it)
# instead of the duplicated "FREQ[23:0] dez" label.
$ws.Range("C2").Value = "FREQ[23:0] hex"

# Move the active selection to A8 (reflects user's next click after editing)
$ws.Range("A8").Select()
